$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# B3 ("Date") must hold the literal text "10.10.2019" (not an actual
# date serial) while keeping its pre-existing date-format style.
# Typing a dd.mm.yyyy-looking string straight into Value gets parsed
# into a real date, so we stash the cell's current formatting on a
# scratch cell, force Text format long enough to type the literal
# string, then restore the saved formatting back onto B3.
# ------------------------------------------------------------------
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "10.10.2019"

$ws.Range("A30").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws.Range("A30").ClearContents() | Out-Null
$ws.Range("A30").ClearFormats() | Out-Null

# Team name / team size
$ws.Range("B4").Value = "MSR Voice Input"
$ws.Range("B5").Value = 5

# Team member names (replacing the generic "Member N" placeholders)
$ws.Range("A8").Value = "Kunaal Sikka"
$ws.Range("A9").Value = "Mina Huh"
$ws.Range("A10").Value = "Vu Nguyen"
$ws.Range("A11").Value = "Nicolas Carmody"
$ws.Range("A12").Value = "Jonas Bokstaller"

# Each member's salary share
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 100

# Tasks completed this week / tasks to complete next week
$ws.Range("A19").Value = "Affinity Clustering"
$ws.Range("B19").Value = "Presentation"
$ws.Range("A20").Value = "Needfinding (mostly)"
$ws.Range("A21").Value = "HoloLens hands-on experience"
$ws.Range("A22").Value = "In person meeting with Stakeholder"
